$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "is_active" column (E, rows 2-7) previously held a =TRUE() boolean
# formula. Replace each with the literal text string "TRUE" (stored as a
# shared string, not a boolean) to fix the boolean-values bug.
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Formula = '="TRUE"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = $false

# Update the active selection/cell on the sheet view to G6.
$ws.Range("G6").Select()
